$wb = $excel.ActiveWorkbook

# --- "작성자명" sheet (Sheet1): fill in the Oct 23 time-log row, then
#     move the selection to F9 (matches the saved selection recorded in
#     the workbook).
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("A10").Value = "10월 23일"
$ws1.Range("B10").Value = 0.5625
$ws1.Range("C10").Value = 0.70833333333333337
$ws1.Range("D10").Value = 0
$ws1.Range("E10").Value = 210
$ws1.Range("F10").Value = "SRS 최종 수정"

# Match the mixed-run formatting used elsewhere in this column: the
# leading "SRS 최종" keeps the cell's base font, " 수정" switches to 돋움 10pt.
$trailingRun = $ws1.Range("F10").Characters(7, 3)
$trailingRun.Font.Name = "돋움"
$trailingRun.Font.Size = 10

[void]$ws1.Range("F9").Select()

# --- Sheet10: set up the page for printing (paper size / orientation).
$ws10 = $wb.Worksheets.Item("Sheet10")
$ws10.PageSetup.PaperSize = 9
$ws10.PageSetup.Orientation = 1
